$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.716.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.290.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.35%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.75%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.288.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("E10").Value = "  +3.65%  "

$ws.Range("E11").Value = "  +3.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.56"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.835.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.53%  "

$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.294.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.769.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("E22").Value = "  +5.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.52%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +4.02%  "

$ws.Range("E31").Value = "  +4.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.14%  "

$ws.Range("E33").Value = "  -2.19%  "

$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("E35").Value = "  +3.26%  "

$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0744"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "425.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.049.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.43%  "

$ws.Range("E42").Value = "  +2.50%  "

$ws.Range("E43").Value = "  +3.26%  "

$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("E46").Value = "  +4.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.95%  "

$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.48%  "

$ws.Range("E51").Value = "  +2.57%  "
